{"js": "// Replace each three-digit-division-by-one-digit answer cell's text with its\n// updated value. Every text run in the document's table cells is targeted;\n// we look each one up with a precise (match-case, whole-text) search and\n// replace it in place so formatting/run properties are preserved.\nconst replacements = [\n  [\"861\u00f72=430, 1\", \"849\u00f78=106, 1\"],\n  [\"339\u00f74=84, 3\", \"144\u00f76=24, 0\"],\n  [\"889\u00f76=148, 1\", \"736\u00f74=184, 0\"],\n  [\"820\u00f77=117, 1\", \"948\u00f72=474, 0\"],\n  [\"476\u00f74=119, 0\", \"963\u00f73=321, 0\"],\n  [\"267\u00f72=133, 1\", \"926\u00f76=154, 2\"],\n  [\"977\u00f74=244, 1\", \"357\u00f77=51, 0\"],\n  [\"407\u00f73=135, 2\", \"536\u00f72=268, 0\"],\n  [\"423\u00f77=60, 3\", \"153\u00f74=38, 1\"],\n  [\"165\u00f79=18, 3\", \"280\u00f77=40, 0\"],\n  [\"483\u00f72=241, 1\", \"287\u00f79=31, 8\"],\n  [\"585\u00f78=73, 1\", \"639\u00f73=213, 0\"],\n  [\"290\u00f77=41, 3\", \"426\u00f76=71, 0\"],\n  [\"205\u00f76=34, 1\", \"142\u00f78=17, 6\"],\n  [\"472\u00f76=78, 4\", \"560\u00f73=186, 2\"],\n  [\"900\u00f78=112, 4\", \"929\u00f77=132, 5\"],\n  [\"803\u00f77=114, 5\", \"283\u00f74=70, 3\"],\n  [\"783\u00f76=130, 3\", \"908\u00f72=454, 0\"],\n  [\"370\u00f73=123, 1\", \"689\u00f79=76, 5\"],\n  [\"943\u00f74=235, 3\", \"871\u00f78=108, 7\"],\n  [\"450\u00f79=50, 0\", \"367\u00f78=45, 7\"],\n  [\"246\u00f72=123, 0\", \"384\u00f76=64, 0\"],\n  [\"800\u00f74=200, 0\", \"757\u00f74=189, 1\"],\n  [\"579\u00f77=82, 5\", \"318\u00f75=63, 3\"],\n  [\"890\u00f73=296, 2\", \"500\u00f75=100, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-division-by-one-digit answer cell's text with its\n# updated value using Word's Find/Replace (Content.Find.Execute), which keeps\n# the existing run formatting (font, size) intact while swapping the text.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"861\u00f72=430, 1\", \"849\u00f78=106, 1\"),\n  @(\"339\u00f74=84, 3\", \"144\u00f76=24, 0\"),\n  @(\"889\u00f76=148, 1\", \"736\u00f74=184, 0\"),\n  @(\"820\u00f77=117, 1\", \"948\u00f72=474, 0\"),\n  @(\"476\u00f74=119, 0\", \"963\u00f73=321, 0\"),\n  @(\"267\u00f72=133, 1\", \"926\u00f76=154, 2\"),\n  @(\"977\u00f74=244, 1\", \"357\u00f77=51, 0\"),\n  @(\"407\u00f73=135, 2\", \"536\u00f72=268, 0\"),\n  @(\"423\u00f77=60, 3\", \"153\u00f74=38, 1\"),\n  @(\"165\u00f79=18, 3\", \"280\u00f77=40, 0\"),\n  @(\"483\u00f72=241, 1\", \"287\u00f79=31, 8\"),\n  @(\"585\u00f78=73, 1\", \"639\u00f73=213, 0\"),\n  @(\"290\u00f77=41, 3\", \"426\u00f76=71, 0\"),\n  @(\"205\u00f76=34, 1\", \"142\u00f78=17, 6\"),\n  @(\"472\u00f76=78, 4\", \"560\u00f73=186, 2\"),\n  @(\"900\u00f78=112, 4\", \"929\u00f77=132, 5\"),\n  @(\"803\u00f77=114, 5\", \"283\u00f74=70, 3\"),\n  @(\"783\u00f76=130, 3\", \"908\u00f72=454, 0\"),\n  @(\"370\u00f73=123, 1\", \"689\u00f79=76, 5\"),\n  @(\"943\u00f74=235, 3\", \"871\u00f78=108, 7\"),\n  @(\"450\u00f79=50, 0\", \"367\u00f78=45, 7\"),\n  @(\"246\u00f72=123, 0\", \"384\u00f76=64, 0\"),\n  @(\"800\u00f74=200, 0\", \"757\u00f74=189, 1\"),\n  @(\"579\u00f77=82, 5\", \"318\u00f75=63, 3\"),\n  @(\"890\u00f73=296, 2\", \"500\u00f75=100, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
